$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("first_eval")

$data = @{
    2  = @(0.2520848399290218, 0.7424843809194045, 1.903216275652063, 1.379571047699995, 1.369577110019525)
    3  = @(0.6974398184037287, 1.552624331747446, 7.423039619458927, 2.72452557695077, 2.659952944235514)
    4  = @(1.179297011432999, 1.950015878589636, 10.42997277502151, 3.22954683740947, 3.037055390127207)
    5  = @(0.9186370495837128, 1.942955877847222, 11.15520434258843, 3.339940769323377, 3.244399679230311)
    6  = @(1.01406298738989, 1.964591607092075, 11.42823783535527, 3.380567679451969, 3.259016548116396)
    7  = @(0.8734069636221409, 1.740984057691856, 10.78249254741741, 3.283670590576559, 3.206762837755778)
    8  = @(1.003986551470107, 1.777346341468097, 11.03323908461691, 3.321631991147862, 3.208769669857039)
    9  = @(0.8645380383347879, 1.931740234399347, 17.21370736197478, 4.148940510777996, 4.158075926358173)
    10 = @(0.03960286749493035, 1.355899730545265, 4.943417352971498, 2.223379714077534, 2.30694412459439)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
}

$wb.Save()
